$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the updated "Price" strings in column D (e.g. "234.71", "0.5993")
# parse as plain numbers, so assigning them straight to Range.Value would let
# Excel auto-convert the cell to a Number (dropping significant trailing zeros
# and changing the stored type away from the original text). To keep those
# cells text, exactly like the source data, we stage each such value in a
# scratch cell that has been explicitly formatted as Text, copy it, and paste
# (values-only) into the destination cell; the destination keeps its own
# original (default) style throughout. The scratch row is removed afterward.
$scratch = $ws.Range("A52")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = '29.095.70'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.823.25'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("E4").Value = '  +0.15%  '
$scratch.Value = '234.71'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -1.83%  '
$scratch.Value = '0.5993'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -3.97%  '
$ws.Range("E7").Value = '  +0.16%  '
$scratch.Value = '0.06951'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -5.75%  '
$scratch.Value = '0.2751'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -4.55%  '
$scratch.Value = '23.25'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -5.94%  '
$scratch.Value = '0.07598'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = '1.837.61'
$ws.Range("E12").Value = '  -0.03%  '
$scratch.Value = '4.738'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -4.20%  '
$scratch.Value = '0.6260'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -5.40%  '
$scratch.Value = '0.000009836'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -6.51%  '
$scratch.Value = '77.33'
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -4.81%  '
$ws.Range("D17").Value = '28.804.83'
$ws.Range("E17").Value = '  -1.64%  '
$scratch.Value = '5.583'
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -10.56%  '
$scratch.Value = '216.19'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -7.59%  '
$scratch.Value = '1.004'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +0.22%  '
$scratch.Value = '11.53'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -5.73%  '
$scratch.Value = '6.873'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -5.74%  '
$scratch.Value = '1.003'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -0.12%  '
$scratch.Value = '156.56'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -0.45%  '
$scratch.Value = '7.929'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -5.73%  '
$scratch.Value = '0.1288'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -3.41%  '
$scratch.Value = '16.45'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -4.55%  '
$scratch.Value = '0.06434'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -9.78%  '
$scratch.Value = '1.424'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -4.03%  '
$scratch.Value = '1.439'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -2.76%  '
$scratch.Value = '3.830'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -4.66%  '
$scratch.Value = '3.754'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -6.67%  '
$scratch.Value = '1.729'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -3.41%  '
$scratch.Value = '1.088'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -5.32%  '
$scratch.Value = '0.6476'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -6.84%  '
$scratch.Value = '2.532'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -2.13%  '
$scratch.Value = '2.737'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -1.69%  '
$scratch.Value = '0.01745'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -4.31%  '
$scratch.Value = '6.529'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -3.70%  '
$ws.Range("D40").Value = '1.145.90'
$ws.Range("E40").Value = '  -7.09%  '
$scratch.Value = '0.8828'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -6.84%  '
$scratch.Value = '1.002'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '1.977.62'
$ws.Range("E43").Value = '  -0.70%  '
$scratch.Value = '100.43'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -0.73%  '
$scratch.Value = '61.60'
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -5.50%  '
$scratch.Value = '0.00000000115'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -1.17%  '
$scratch.Value = '1.609'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -4.34%  '
$scratch.Value = '8.490'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -4.69%  '
$ws.Range("E49").Value = '  -2.53%  '
$scratch.Value = '0.4531'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -0.88%  '
$scratch.Value = '6.419'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -7.29%  '

$excel.CutCopyMode = $false
$scratch.EntireRow.Delete()
